# Horarios Linea 141 - actualizacion de datos (scrap 17:54:43)
# Actualiza las 3 hojas (LP1912, LP1912-215, 6203-6173) con las nuevas
# filas raspadas: se insertan nuevas llegadas (ordenadas por Hora_Llegada),
# lo que desplaza hacia abajo algunas filas existentes y agrega filas nuevas
# al final de cada tabla. Tambien se actualizan los encabezados
# "Ultima actualizacion" y "Total filas".

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Hoja "LP1912"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LP1912")

$ws.Cells.Item(2,1).Value = "Última actualización: 17:54:43"
$ws.Cells.Item(3,1).Value = "Total filas: 72"

$ws.Cells.Item(55,1).Value = "17:54:43"
$ws.Cells.Item(55,2).Value = "18:33"
$ws.Cells.Item(55,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(55,4).Value = 39
$ws.Cells.Item(55,5).Value = "LP1912"

$ws.Cells.Item(56,1).Value = "16:46:42"
$ws.Cells.Item(56,2).Value = "18:34"
$ws.Cells.Item(56,3).Value = "14X44_ABASTO"
$ws.Cells.Item(56,4).Value = 108
$ws.Cells.Item(56,5).Value = "LP1912"

$ws.Cells.Item(57,1).Value = "16:46:42"
$ws.Cells.Item(57,2).Value = "18:38"
$ws.Cells.Item(57,3).Value = "17X38_ROMERO"
$ws.Cells.Item(57,4).Value = 112
$ws.Cells.Item(57,5).Value = "LP1912"

$ws.Cells.Item(58,1).Value = "17:13:30"
$ws.Cells.Item(58,2).Value = "18:41"
$ws.Cells.Item(58,3).Value = "14_ABASTO"
$ws.Cells.Item(58,4).Value = 88
$ws.Cells.Item(58,5).Value = "LP1912"

$ws.Cells.Item(59,1).Value = "16:46:42"
$ws.Cells.Item(59,2).Value = "18:41"
$ws.Cells.Item(59,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(59,4).Value = 115
$ws.Cells.Item(59,5).Value = "LP1912"

$ws.Cells.Item(60,1).Value = "17:47:22"
$ws.Cells.Item(60,2).Value = "18:44"
$ws.Cells.Item(60,3).Value = "14_ABASTO"
$ws.Cells.Item(60,4).Value = 57
$ws.Cells.Item(60,5).Value = "LP1912"

$ws.Cells.Item(61,1).Value = "17:35:09"
$ws.Cells.Item(61,2).Value = "18:45"
$ws.Cells.Item(61,3).Value = "14_ABASTO"
$ws.Cells.Item(61,4).Value = 70
$ws.Cells.Item(61,5).Value = "LP1912"

$ws.Cells.Item(62,1).Value = "17:35:09"
$ws.Cells.Item(62,2).Value = "18:51"
$ws.Cells.Item(62,3).Value = "15_ABASTO"
$ws.Cells.Item(62,4).Value = 76
$ws.Cells.Item(62,5).Value = "LP1912"

$ws.Cells.Item(63,1).Value = "17:54:43"
$ws.Cells.Item(63,2).Value = "18:53"
$ws.Cells.Item(63,3).Value = "16_SANTA ANA"
$ws.Cells.Item(63,4).Value = 59
$ws.Cells.Item(63,5).Value = "LP1912"

$ws.Cells.Item(64,1).Value = "17:35:09"
$ws.Cells.Item(64,2).Value = "18:59"
$ws.Cells.Item(64,3).Value = "10_OLMOS"
$ws.Cells.Item(64,4).Value = 84
$ws.Cells.Item(64,5).Value = "LP1912"

$ws.Cells.Item(65,1).Value = "17:13:30"
$ws.Cells.Item(65,2).Value = "19:01"
$ws.Cells.Item(65,3).Value = "17_ROMERO"
$ws.Cells.Item(65,4).Value = 108
$ws.Cells.Item(65,5).Value = "LP1912"

$ws.Cells.Item(66,1).Value = "17:13:30"
$ws.Cells.Item(66,2).Value = "19:11"
$ws.Cells.Item(66,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(66,4).Value = 118
$ws.Cells.Item(66,5).Value = "LP1912"

$ws.Cells.Item(67,1).Value = "17:47:22"
$ws.Cells.Item(67,2).Value = "19:17"
$ws.Cells.Item(67,3).Value = "27_EL RETIRO"
$ws.Cells.Item(67,4).Value = 90
$ws.Cells.Item(67,5).Value = "LP1912"

$ws.Cells.Item(68,1).Value = "17:35:09"
$ws.Cells.Item(68,2).Value = "19:19"
$ws.Cells.Item(68,3).Value = "27_EL RETIRO"
$ws.Cells.Item(68,4).Value = 104
$ws.Cells.Item(68,5).Value = "LP1912"

$ws.Cells.Item(69,1).Value = "17:54:43"
$ws.Cells.Item(69,2).Value = "19:20"
$ws.Cells.Item(69,3).Value = "215C_EL PATO"
$ws.Cells.Item(69,4).Value = 86
$ws.Cells.Item(69,5).Value = "LP1912"

$ws.Cells.Item(70,1).Value = "17:35:09"
$ws.Cells.Item(70,2).Value = "19:21"
$ws.Cells.Item(70,3).Value = "215C_EL PATO"
$ws.Cells.Item(70,4).Value = 106
$ws.Cells.Item(70,5).Value = "LP1912"

$ws.Cells.Item(71,1).Value = "17:35:09"
$ws.Cells.Item(71,2).Value = "19:29"
$ws.Cells.Item(71,3).Value = "225_GOMEZ"
$ws.Cells.Item(71,4).Value = 114
$ws.Cells.Item(71,5).Value = "LP1912"

$ws.Cells.Item(72,1).Value = "17:54:43"
$ws.Cells.Item(72,2).Value = "19:30"
$ws.Cells.Item(72,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(72,4).Value = 96
$ws.Cells.Item(72,5).Value = "LP1912"

$ws.Cells.Item(73,1).Value = "17:54:43"
$ws.Cells.Item(73,2).Value = "19:30"
$ws.Cells.Item(73,3).Value = "27_EL RETIRO"
$ws.Cells.Item(73,4).Value = 96
$ws.Cells.Item(73,5).Value = "LP1912"

$ws.Cells.Item(74,1).Value = "17:35:09"
$ws.Cells.Item(74,2).Value = "19:31"
$ws.Cells.Item(74,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(74,4).Value = 116
$ws.Cells.Item(74,5).Value = "LP1912"

$ws.Cells.Item(75,1).Value = "17:47:22"
$ws.Cells.Item(75,2).Value = "19:40"
$ws.Cells.Item(75,3).Value = "17X38_ROMERO"
$ws.Cells.Item(75,4).Value = 113
$ws.Cells.Item(75,5).Value = "LP1912"

$ws.Cells.Item(76,1).Value = "17:47:22"
$ws.Cells.Item(76,2).Value = "19:44"
$ws.Cells.Item(76,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(76,4).Value = 117
$ws.Cells.Item(76,5).Value = "LP1912"

$ws.Cells.Item(77,1).Value = "17:54:43"
$ws.Cells.Item(77,2).Value = "19:51"
$ws.Cells.Item(77,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(77,4).Value = 117
$ws.Cells.Item(77,5).Value = "LP1912"

# ------------------------------------------------------------------
# Hoja "LP1912-215"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Cells.Item(2,1).Value = "Última actualización: 17:54:43"
$ws.Cells.Item(3,1).Value = "Total filas: 10"

$ws.Cells.Item(12,1).Value = "17:54:43"
$ws.Cells.Item(12,2).Value = "19:20"
$ws.Cells.Item(12,3).Value = "215C_EL PATO"
$ws.Cells.Item(12,4).Value = 86
$ws.Cells.Item(12,5).Value = "LP1912"

$ws.Cells.Item(13,1).Value = "17:35:09"
$ws.Cells.Item(13,2).Value = "19:21"
$ws.Cells.Item(13,3).Value = "215C_EL PATO"
$ws.Cells.Item(13,4).Value = 106
$ws.Cells.Item(13,5).Value = "LP1912"

$ws.Cells.Item(14,1).Value = "17:54:43"
$ws.Cells.Item(14,2).Value = "19:30"
$ws.Cells.Item(14,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(14,4).Value = 96
$ws.Cells.Item(14,5).Value = "LP1912"

$ws.Cells.Item(15,1).Value = "17:35:09"
$ws.Cells.Item(15,2).Value = "19:31"
$ws.Cells.Item(15,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(15,4).Value = 116
$ws.Cells.Item(15,5).Value = "LP1912"

# ------------------------------------------------------------------
# Hoja "6203-6173"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("6203-6173")

$ws.Cells.Item(2,1).Value = "Última actualización: 17:54:43"
